$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.087.62"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "3.849.22"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "698.43"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.74"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("D7").Value = "3.847.37"
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.36"
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.84"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "4.499.22"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "3.875.62"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "71.230.77"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.25"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.47"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "498.09"
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.68"
$ws.Range("E22").Value = "  -3.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.738"
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.39"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000145"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.65"
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.18"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.10"
$ws.Range("E28").Value = "  -3.31%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.07"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.47"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.24"
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.43"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.177"
$ws.Range("E34").Value = "  -5.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.22"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("D36").Value = "3.813.51"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("E40").Value = "  +4.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.01"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.36"
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("E45").Value = "  +3.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.65"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "432.97"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.03"
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.74"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.298"
$ws.Range("E51").Value = "  -1.97%  "
